$wb = $excel.ActiveWorkbook

# Column F values (想去人数 / "want to go" count) changed for rows 3-6
# in both the "展览" and "全部类型" sheets.
$updates = @{
    "F3" = 1276
    "F4" = 1584
    "F5" = 63
    "F6" = 6191
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
